# manuallyEnteredGroceryReceipts.xlsx — classify June transactions (rows 47-101)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# item, amount, date-serial for each new receipt line (rows 47 through 101)
$newRows = @(
    @(47, "Bread", 8.4, 43628),
    @(48, "Chocolate", 5.3, 43627),
    @(49, "Nappies", 15, 43627),
    @(50, "Milk", 3.41, 43627),
    @(51, "Milk", 3.41, 43627),
    @(52, "Laundry liquid", 7.99, 43627),
    @(53, "Only organic lasagne", 5.4, 43627),
    @(54, "Bathroom bags", 3.6, 43627),
    @(55, "Free range chicken", 7.99, 43627),
    @(56, "Only organic blueberry rice cakes", 2.99, 43627),
    @(57, "Bread", 4.2, 43626),
    @(58, "Moisturiser", 36, 43626),
    @(59, "Nappies", 7, 43625),
    @(60, "Reusable bag", 1, 43625),
    @(61, "Popsicles", 4.9000000000000004, 43625),
    @(62, "Bagels", 4.3, 43625),
    @(63, "Nappies", 10, 43625),
    @(64, "Apple juice", 1.7, 43625),
    @(65, "Soup", 4.5, 43625),
    @(66, "Tomatoes", 1.4, 43625),
    @(67, "Tomato paste", 4.5, 43625),
    @(68, "Red onions", 2.7, 43625),
    @(69, "Garlic", 4, 43625),
    @(70, "Free range eggs", 7, 43625),
    @(71, "Chicken thighs", 7.9, 43625),
    @(72, "Chicken thighs", 6.7, 43625),
    @(73, "Chicken thighs", 8.9, 43625),
    @(74, "lentils", 6.82, 43624),
    @(75, "Turtle beans", 7.25, 43624),
    @(76, "Kidney beans", 6.65, 43624),
    @(77, "Paprika", 3.06, 43624),
    @(78, "Whole peppercorns", 3.37, 43624),
    @(79, "Anchor dairy blend", 4.5999999999999996, 43623),
    @(80, "Anchor cheese", 9.5, 43623),
    @(81, "Hairspray", 9, 43623),
    @(82, "Nappies", 10, 43620),
    @(83, "Tomatoes", 4.2, 43620),
    @(84, "Tomato paste", 4.49, 43620),
    @(85, "Tomato paste", 4.49, 43620),
    @(86, "Peanut butter", 6, 43620),
    @(87, "Beef mince", 9, 43620),
    @(88, "Brown onions", 5.98, 43620),
    @(89, "Garlic", 4, 43620),
    @(90, "kumara", 1.81, 43620),
    @(91, "Bread", 8.4, 43619),
    @(92, "Pita bread", 4.5999999999999996, 43619),
    @(93, "Avalanche coffee", 6, 43619),
    @(94, "Banana porridge", 3, 43619),
    @(95, "Only organic apple and cinnamon biscotti", 4, 43619),
    @(96, "Chocolate", 4.5, 43619),
    @(97, "Licorice allsorts", 4, 43619),
    @(98, "Only organic blueberry rice cake", 3.69, 43619),
    @(99, "Panadol", 4, 43619),
    @(100, "Dental floss", 3.5, 43619),
    @(101, "Corn chips", 3.6, 43619)
)

foreach ($r in $newRows) {
    $rowNum = $r[0]
    $ws.Range("A$rowNum").Value = $r[1]
    $ws.Range("B$rowNum").Value = $r[2]
    $ws.Range("C$rowNum").Value = $r[3]
}

# --- Match existing formatting conventions for the new rows ---
# Column A (item): same style as the other item cells (e.g. A4)
$ws.Range("A4").Copy() | Out-Null
$ws.Range("A47:A101").PasteSpecial(-4122) | Out-Null
# Row 56 item cell uses the wrap-text style (like A3)
$ws.Range("A3").Copy() | Out-Null
$ws.Range("A56").PasteSpecial(-4122) | Out-Null

# Column C (date): date-format style matching the other date cells (e.g. C3)
$ws.Range("C3").Copy() | Out-Null
$ws.Range("C47:C101").PasteSpecial(-4122) | Out-Null

# Column B (amount): rows 47 and 79 carry the styled-amount format (like B3)
$ws.Range("B3").Copy() | Out-Null
$ws.Range("B47").PasteSpecial(-4122) | Out-Null
$ws.Range("B79").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = $false

# --- Update the view: scroll position and active selection ---
$win = $excel.ActiveWindow
$win.ScrollRow = 72
$win.ScrollColumn = 1
$ws.Range("C92").Select() | Out-Null
